# Auto-applied numeric updates from scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1022176.06
$ws.Range("I15").Value = 1022176.06
$ws.Range("K15").Value = 3066528.18
$ws.Range("M15").Value = -3066359.18

$ws.Range("H40").Value = 3662.3635
$ws.Range("I40").Value = 4037.75
$ws.Range("J40").Value = 3447.8572
$ws.Range("K40").Value = 4037.75
$ws.Range("L40").Value = 3447.8572
$ws.Range("M40").Value = -3862.75
$ws.Range("N40").Value = -3797.8572

$ws.Range("H51").Value = 15155513
$ws.Range("I51").Value = 30306362
$ws.Range("J51").Value = 4663.3335
$ws.Range("K51").Value = 30306362
$ws.Range("L51").Value = 4663.3335
$ws.Range("M51").Value = -30305878
$ws.Range("N51").Value = -5631.3335

$ws.Range("H64").Value = 3045.6191
$ws.Range("I64").Value = 2978
$ws.Range("J64").Value = 3096.3333
$ws.Range("K64").Value = 2978
$ws.Range("L64").Value = 3096.3333
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -3592.3333

$ws.Range("H67").Value = 3045.6191
$ws.Range("I67").Value = 2978
$ws.Range("J67").Value = 3096.3333
$ws.Range("K67").Value = 2978
$ws.Range("L67").Value = 3096.3333
$ws.Range("M67").Value = -2120
$ws.Range("N67").Value = -4812.3333

$ws.Range("H76").Value = 2894.5144
$ws.Range("I76").Value = 2752.4075
$ws.Range("J76").Value = 3374.125
$ws.Range("K76").Value = 2752.4075
$ws.Range("L76").Value = 3374.125
$ws.Range("M76").Value = -2437.4075
$ws.Range("N76").Value = -4004.125

$ws.Range("H79").Value = 2894.5144
$ws.Range("I79").Value = 2752.4075
$ws.Range("J79").Value = 3374.125
$ws.Range("K79").Value = 2752.4075
$ws.Range("L79").Value = 3374.125
$ws.Range("M79").Value = -1660.4075
$ws.Range("N79").Value = -5558.125

$ws.Range("H112").Value = 2558.925
$ws.Range("J112").Value = 2688.027
$ws.Range("L112").Value = 8064.081
$ws.Range("N112").Value = -10280.081

$ws.Range("H116").Value = 2382.8696
$ws.Range("I116").Value = 1877.7778
$ws.Range("K116").Value = 1877.7778
$ws.Range("M116").Value = 1564.2222

$ws.Range("H137").Value = 1236.9077
$ws.Range("I137").Value = 1797.4814
$ws.Range("J137").Value = 838.6053
$ws.Range("K137").Value = 5392.4442
$ws.Range("L137").Value = 2515.8159
$ws.Range("M137").Value = -2842.4442
$ws.Range("N137").Value = -7615.8159

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1328.6666
$ws.Range("I16").Value = 1328.6666
$ws.Range("K16").Value = 1328.6666
$ws.Range("M16").Value = -1041.6666

$ws.Range("H32").Value = 4967.18
$ws.Range("I32").Value = 3974.5212
$ws.Range("J32").Value = 20518.834
$ws.Range("K32").Value = 3974.5212
$ws.Range("L32").Value = 20518.834
$ws.Range("M32").Value = -3687.5212
$ws.Range("N32").Value = -21092.834

$ws.Range("H63").Value = 2901.5715
$ws.Range("I63").Value = 2098.5186
$ws.Range("J63").Value = 7719.8887
$ws.Range("K63").Value = 2098.5186
$ws.Range("L63").Value = 7719.8887
$ws.Range("M63").Value = -1412.5186
$ws.Range("N63").Value = -9091.8887

$ws.Range("H66").Value = 2901.5715
$ws.Range("I66").Value = 2098.5186
$ws.Range("J66").Value = 7719.8887
$ws.Range("K66").Value = 10492.593
$ws.Range("L66").Value = 38599.4435
$ws.Range("M66").Value = -7060.592999999999
$ws.Range("N66").Value = -45463.4435

$ws.Range("H113").Value = 34999
$ws.Range("J113").Value = 34999
$ws.Range("L113").Value = 34999
$ws.Range("N113").Value = -43677

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1885.9048
$ws.Range("I105").Value = 1557.2727
$ws.Range("K105").Value = 1557.2727
$ws.Range("M105").Value = 189.7273

$ws.Range("H139").Value = 64983.332
$ws.Range("J139").Value = 64983.332
$ws.Range("L139").Value = 64983.332
$ws.Range("N139").Value = -75263.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2312.6562
$ws.Range("I31").Value = 1423.3055
$ws.Range("J31").Value = 3456.1072
$ws.Range("K31").Value = 1423.3055
$ws.Range("L31").Value = 3456.1072
$ws.Range("M31").Value = -1128.3055
$ws.Range("N31").Value = -4046.1072

$ws.Range("H34").Value = 2312.6562
$ws.Range("I34").Value = 1423.3055
$ws.Range("J34").Value = 3456.1072
$ws.Range("K34").Value = 1423.3055
$ws.Range("L34").Value = 3456.1072
$ws.Range("M34").Value = -1221.3055
$ws.Range("N34").Value = -3860.1072

$ws.Range("H62").Value = 6982.9565
$ws.Range("I62").Value = 2598.25
$ws.Range("J62").Value = 17005.143
$ws.Range("K62").Value = 2598.25
$ws.Range("L62").Value = 17005.143
$ws.Range("M62").Value = -1974.25
$ws.Range("N62").Value = -18253.143

$ws.Range("H65").Value = 6982.9565
$ws.Range("I65").Value = 2598.25
$ws.Range("J65").Value = 17005.143
$ws.Range("K65").Value = 12991.25
$ws.Range("L65").Value = 85025.715
$ws.Range("M65").Value = -9871.25
$ws.Range("N65").Value = -91265.715

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 752.36365
$ws.Range("I132").Value = 636.5357
$ws.Range("J132").Value = 1401
$ws.Range("K132").Value = 5728.821300000001
$ws.Range("L132").Value = 12609
$ws.Range("M132").Value = -3198.821300000001
$ws.Range("N132").Value = -17669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 205002100
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10336

$ws.Range("H70").Value = 7627.6
$ws.Range("I70").Value = 10575.4
$ws.Range("J70").Value = 4679.8
$ws.Range("K70").Value = 10575.4
$ws.Range("L70").Value = 4679.8
$ws.Range("M70").Value = -10305.4
$ws.Range("N70").Value = -5219.8

$ws.Range("H73").Value = 7627.6
$ws.Range("I73").Value = 10575.4
$ws.Range("J73").Value = 4679.8
$ws.Range("K73").Value = 10575.4
$ws.Range("L73").Value = 4679.8
$ws.Range("M73").Value = -9639.4
$ws.Range("N73").Value = -6551.8

$ws.Range("H80").Value = 2542.9375
$ws.Range("I80").Value = 2368.2666
$ws.Range("J80").Value = 2697.0588
$ws.Range("K80").Value = 2368.2666
$ws.Range("L80").Value = 2697.0588
$ws.Range("M80").Value = -1370.2666
$ws.Range("N80").Value = -4693.0588

$ws.Range("H83").Value = 2542.9375
$ws.Range("I83").Value = 2368.2666
$ws.Range("J83").Value = 2697.0588
$ws.Range("K83").Value = 11841.333
$ws.Range("L83").Value = 13485.294
$ws.Range("M83").Value = -6849.332999999999
$ws.Range("N83").Value = -23469.294

$ws.Range("H105").Value = 33835.5
$ws.Range("J105").Value = 33835.5
$ws.Range("L105").Value = 33835.5
$ws.Range("N105").Value = -40823.5

$ws.Range("H110").Value = 44950
$ws.Range("J110").Value = 44950
$ws.Range("L110").Value = 44950
$ws.Range("N110").Value = -53130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 28538.334
$ws.Range("J105").Value = 28538.334
$ws.Range("L105").Value = 28538.334
$ws.Range("N105").Value = -35526.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3531.1428
$ws.Range("I62").Value = 2737.5
$ws.Range("J62").Value = 3848.6
$ws.Range("K62").Value = 2737.5
$ws.Range("L62").Value = 3848.6
$ws.Range("M62").Value = -2113.5
$ws.Range("N62").Value = -5096.6

$ws.Range("H65").Value = 3531.1428
$ws.Range("I65").Value = 2737.5
$ws.Range("J65").Value = 3848.6
$ws.Range("K65").Value = 13687.5
$ws.Range("L65").Value = 19243
$ws.Range("M65").Value = -10567.5
$ws.Range("N65").Value = -25483

$ws.Range("H104").Value = 21656.666
$ws.Range("J104").Value = 21656.666
$ws.Range("L104").Value = 21656.666
$ws.Range("N104").Value = -28644.666
